$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Replace-Text: search text not found: $old"
    }
}

# 1) Combined VISTOS paragraph (Formulario, Oficio, two Proveidos, acuerdo CU)
Replace-Text "El Formulario Único de Trámite N° 0001-CF-FCC.SS.-UPLA-2023 de fecha 30.06.2023, Oficio N° 0002-E-2023-OURYM-UPLA de fecha 15.06.2023, Proveído N° 0003-2023-VRACD-UPLA de fecha 30.06.2023, Proveído N° 0004-2023-VRACD-UPLA de fecha 30.06.2023 y acuerdo de Consejo Universitario en sesión extraordinaria de fecha 30.06.2023, respectivamente; y," "El Formulario Único de Trámite N° 0002-CF-FCC.SS.-UPLA-2023 de fecha 03.07.2023, Oficio N° 0003-E-2023-OURYM-UPLA de fecha 29.06.2023, Proveído N° 0004-2023-VRACD-UPLA de fecha 26.06.2023, Proveído N° 0005-2023-VRACD-UPLA de fecha 25.06.2023 y acuerdo de Consejo Universitario en sesión extraordinaria de fecha 03.07.2023, respectivamente; y,"

# 2) CONSIDERANDO paragraph: Formulario Único de Trámite N° 0001 de fecha 30.06.2023
Replace-Text "El Formulario Único de Trámite N° 0001 de fecha 30.06.2023, presentado por el señor Vargas Cruz Juan" "El Formulario Único de Trámite N° 0002 de fecha 03.07.2023, presentado por el señor Vargas Cruz Juan"

# 3) CONSIDERANDO paragraph: Oficio N° 0002-E-2023-OURYM-UPLA de fecha 15.06.2023
Replace-Text "El Oficio N° 0002-E-2023-OURYM-UPLA de fecha 15.06.2023 emitido por el Jefe" "El Oficio N° 0003-E-2023-OURYM-UPLA de fecha 29.06.2023 emitido por el Jefe"

# 4) CONSIDERANDO paragraph: Los Miembros del Consejo Universitario ... 30.06.2023
Replace-Text "Los Miembros del Consejo Universitario en sesión extraordinaria de fecha 30.06.2023, toman conocimiento" "Los Miembros del Consejo Universitario en sesión extraordinaria de fecha 03.07.2023, toman conocimiento"

# 5) MODALIDAD block (DICE / DEBE DECIR) - the paragraph uses embedded line-feed (Chr(10))
#    characters rather than Word paragraph marks, so build the search/replace strings
#    with an explicit LF rather than the ^p wildcard. (Build into plain variables
#    first -- passing parenthesized expressions directly as call args trips this
#    interpreter up.)
$nl = [char]10
$indent = "                "
$old5 = "DICE:" + $nl + $indent + "MODALIDAD: Primera Selección"
$new5 = "DICE:" + $nl + $indent + "MODALIDAD: Postulante Ordinario o Regular"
Replace-Text $old5 $new5

$old6 = "DEBE DECIR: " + $nl + $indent + "MODALIDAD: Primera Selección"
$new6 = "DEBE DECIR: " + $nl + $indent + "MODALIDAD: Mayores de 30 años"
Replace-Text $old6 $new6

# 6) Header date "Huancayo, 30.06.2023" -> "03.07.2023"
#    $d.Content only covers the main story; the date lives in the primary
#    header of the (single) section, so reach it via Sections/Headers.
$sec = $d.Sections.First
$hdr = $sec.Headers.Item(1)
if ($hdr.Exists) {
    $hdrFound = $hdr.Range.Find.Execute("30.06.2023", $true, $true, $false, $false, $false, $true, 1, $false, "03.07.2023", 2)
    if (-not $hdrFound) {
        throw "Header date replacement: '30.06.2023' not found in primary header"
    }
}
else {
    throw "Primary header does not exist"
}

Write-Output "Done"
